$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of data (6th of May 2020 stats) after the existing last row (row 24)
$ws.Cells.Item(25, 1).Value = 43957
$ws.Cells.Item(25, 2).Value = 82
$ws.Cells.Item(25, 3).Value = 37
$ws.Cells.Item(25, 4).Value = 144
$ws.Cells.Item(25, 5).Value = 52

# Update the active selection to match the newly added row, as Excel would after data entry
$ws.Range("C25").Select()
